# Updated cryptos list values (price + 1h volume change) and
# two pairs of rows whose ranking order swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.164.00"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -6.72%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.671.52"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -4.18%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.23%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "218.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.79%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5068"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -12.49%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.005"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2641"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.67%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06346"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -3.76%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.55"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -7.01%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07370"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.02%  "

$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.563"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -3.40%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.668.01"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -4.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5775"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -4.60%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.891.54"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -4.61%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.000008514"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.03%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.00"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -12.57%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.166.48"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -6.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.958"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -7.03%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.007"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.46%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.81"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -4.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "191.12"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -6.87%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.202"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -6.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.005"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.21%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "142.78"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -4.80%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.691"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -4.65%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1175"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.88"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.56%  "

$ws.Range("B29").Value = "Hedera"
$ws.Range("C29").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05854"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -5.87%  "

$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.295"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -6.64%  "

$ws.Range("E31").Value = "  -5.08%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.508"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.35%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.515"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -5.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.656"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.29%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.009"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.71%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6003"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -5.93%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.361"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.66%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.645"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -3.10%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01605"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.08%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.086.67"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -3.78%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.951"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -4.06%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8598"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.66%  "

$ws.Range("E43").Value = "  +0.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.60"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.815.89"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -4.20%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000111"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +3.85%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.94"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -5.95%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.008"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.92%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.049"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -2.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4294"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.87%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05176"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -3.74%  "
